$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric values (e.g. "1.000", "0.9999") are kept as
# literal text rather than being auto-converted to numbers by Excel,
# matching the inlineStr cells in the source workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.350.23"
$ws.Range("E2").Value = "  -3.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.939.29"
$ws.Range("E3").Value = "  -3.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.85"
$ws.Range("E5").Value = "  -3.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7230"
$ws.Range("E6").Value = "  -6.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3351"
$ws.Range("E8").Value = "  -5.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.56"
$ws.Range("E9").Value = "  -0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07334"
$ws.Range("E10").Value = "  +3.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08130"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.37"
$ws.Range("E15").Value = "  -6.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.93"
$ws.Range("E16").Value = "  -4.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.340.18"
$ws.Range("E17").Value = "  -3.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008323"
$ws.Range("E18").Value = "  +4.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.92"
$ws.Range("E19").Value = "  -7.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.872"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.188.46"
$ws.Range("E21").Value = "  -3.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.957"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.844"
$ws.Range("E25").Value = "  -2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.63"
$ws.Range("E26").Value = "  -2.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.423"
$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.73"
$ws.Range("E28").Value = "  -1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1333"
$ws.Range("E29").Value = "  -9.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.562"
$ws.Range("E30").Value = "  -3.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.451"
$ws.Range("E32").Value = "  -4.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.249"
$ws.Range("E33").Value = "  -4.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05236"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.267"
$ws.Range("E35").Value = "  +3.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7533"
$ws.Range("E36").Value = "  -3.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.740"
$ws.Range("E37").Value = "  -2.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01999"
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.844"
$ws.Range("E39").Value = "  -3.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.680"
$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "79.99"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4554"
$ws.Range("E42").Value = "  -3.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.035"
$ws.Range("E43").Value = "  -5.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.96"
$ws.Range("E46").Value = "  -4.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.815"
$ws.Range("E47").Value = "  -1.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.424"
$ws.Range("E48").Value = "  -4.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.93"
$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.504"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4157"
$ws.Range("E51").Value = "  -4.26%  "

$ws.Range("E11").Value = "  -5.13%  "
$ws.Range("E31").Value = "  -1.28%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.936.82"
$ws.Range("E13").Value = "  -3.52%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.531"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8429"
$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.31%  "
